# "squash should be rebase"
#
# Slide 16, "Content Placeholder 2" shape, 2nd paragraph currently reads:
#   "Git squash – squashes all the changes in a branch down to one."
# It needs to become two runs reading:
#   "Git rebase – squashes all the changes in a branch down to one, or rebases on "
#   "a different base"

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(16)
$shape = $slide.Shapes.Item(2)
$textRange = $shape.TextFrame.TextRange

$targetParagraph = $textRange.Paragraphs(2)

# Rewrite the existing run's text (keeps its original run properties,
# e.g. lang="en-US" dirty="0").
$firstRun = $targetParagraph.Runs(1)
$firstRun.Text = "Git rebase " + [char]0x2013 + " squashes all the changes in a branch down to one, or rebases on "

# Append the new trailing run with the additional clause.
$secondRun = $targetParagraph.InsertAfter("a different base")

Write-Host "Updated paragraph text:" $targetParagraph.Text
